$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Security Groups sheet: append the new Production rows (PR 3379 / WI #31936)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Security Groups")

$newRows = @(
    @("Production", "Administrators", "WFM\IRMA.Deploy.Prod",      "Job Web"),
    @("Production", "Administrators", "WFM\Icon.Deploy.Prod",      "Job Web"),
    @("Production", "Administrators", "WFM\Mammoth.Deploy.Prod",   "Job Web"),
    @("Production", "Administrators", "wfm\IconInterfaceUserPrd",  "Job Web"),
    @("Production", "Administrators", "wfm\MammothPrd",            "Job Web"),
    @("Production", "Administrators", "WFM\SPOReportsDev",         "Web"),
    @("Production", "Administrators", "WFM\NutriconService",       "Web"),
    @("Production", "Administrators", "wfm\IconWebPrd",            "Web"),
    @("Production", "IIS_IUSRS",      "wfm\MammothPrd",            "Web"),
    @("Production", "IIS_IUSRS",      "wfm\IconWebPrd",            "Web"),
    @("Production", "IIS_IUSRS",      "WFM\NutriconService",       "Web"),
    @("Production", "IIS_IUSRS",      "WFM\SPOReports",            "Web"),
    @("Production", "IIS_IUSRS",      "Authenticated Users",       "Web")
)

$startRow = 44
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
}

# New column-A width (typed values in col A widened it beyond its default)
$ws.Columns.Item(1).ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# View-state: user's last interaction left the selection on "Security
# Groups"!C55 (scrolled to row 25) and "Server List"!G72.
# ---------------------------------------------------------------------------
$wsServers = $wb.Worksheets.Item("Server List")
$wsServers.Activate()
$wsServers.Range("G72").Select()

$ws.Activate()
$ws.Range("C55").Select()

Write-Output "done"
